# Update "想去人数" (F column) counts for matching event rows on both the
# "展览" sheet and the "全部类型" sheet, per the data refresh commit.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> row number -> new F value
$updates = @{
    "展览" = @{
        4  = 3682
        5  = 2250
        6  = 438
        12 = 1364
        14 = 2086
        15 = 150
    }
    "全部类型" = @{
        4  = 3682
        5  = 2250
        6  = 438
        15 = 1364
        17 = 2086
        18 = 150
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $ws.Cells.Item($rowNum, 6).Value = $rows[$rowNum]
    }
}
